$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- J14:J16 previously-empty accuracy values now populated ---
$ws.Range("J14").Value = 0.989
$ws.Range("J15").Value = 0.879
$ws.Range("J16").Value = 0.877

# --- K14:M16 timing formulas updated from placeholder 0/60 to real totals ---
$ws.Range("K14").Formula = "=734.58/60"
$ws.Range("L14").Formula = "=15406.02/60"
$ws.Range("M14").Formula = "=28.26/60"

$ws.Range("K15").Formula = "=1017.62/60"
$ws.Range("L15").Formula = "=21060.58/60"
$ws.Range("M15").Formula = "=34.87/60"

$ws.Range("K16").Formula = "=1339.89/60"
$ws.Range("L16").Formula = "=39843.97/60"
$ws.Range("M16").Formula = "=47.47/60"

# --- Row 17 (520000-row model) fully populated to match rows 13:16 above it ---
$ws.Range("B17").Formula = "=A17-D17"
$ws.Range("C17").Formula = "=B17/A17"
$ws.Range("D17").Value = 371670
$ws.Range("E17").Value = 0.908
$ws.Range("F17").Value = 0.927
$ws.Range("G17").Value = 0.927
$ws.Range("H17").Value = 0.535
$ws.Range("I17").Value = 0.909
$ws.Range("J17").Value = 0.874
$ws.Range("K17").Formula = "=1672.72/60"
$ws.Range("L17").Formula = "=33395.93/60"
$ws.Range("M17").Formula = "=51.26/60"
$ws.Range("N17").Value = 130

# Match the number-format style used by the rest of the "%NAs" column (C13:C16)
$ws.Range("C17").NumberFormat = $ws.Range("C16").NumberFormat

# --- View state: scrolled down a couple rows, selection moved from A14 to A16 ---
$ws.Range("A16").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 11
$activeWindow.ScrollColumn = 1
